$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 728, shifting existing rows 728:840 down to 729:841
$ws.Rows.Item(728).Insert()

# Populate the newly inserted row 728 with the new record
$ws.Range("A728").Value = 4
$ws.Range("B728").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C728").Value = "Los Lagos"
$ws.Range("D728").Value = 45218
$ws.Range("E728").Value = 10
$ws.Range("F728").Value = 100112006
$ws.Range("G728").Value = "Repollo"
$ws.Range("H728").Value = "Crespo record"
$ws.Range("I728").Value = "Primera"
$ws.Range("J728").Value = 500
$ws.Range("K728").Value = 1500
$ws.Range("L728").Value = 1500
$ws.Range("M728").Value = 1500
$ws.Range("N728").Value = "$/unidad"
$ws.Range("O728").Value = "Región Metropolitana"
$ws.Range("P728").Value = 1500
$ws.Range("Q728").Value = 1
$ws.Range("R728").Value = "Hortaliza"
